# "Update the pres y0"
#
# 1. Insert a new slide ("Our approach cont..") right before the
#    "Alternatives" slide (i.e. as the new slide 6).
# 2. Add a trailing full stop to the last three bullets of the
#    "Preferred Content" slide (slide 4).
# 3. Re-cache the "today" date field shown on the master/layouts
#    (01/02/2013 -> 06/02/2013).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. New slide, inserted at position 6 (before "Alternatives"),
#    using the "Title and Content" layout (same as the other
#    bullet-point slides in the deck).
# ---------------------------------------------------------------
$newSlide = $p.Slides.Add(6, 2)

$title = $newSlide.Shapes.Item(2)
$title.TextFrame.TextRange.Text = "Our approach cont.."

$bodyLines = @(
    "",
    "For 3D renders, we will be using blender as it is open-source and works with the XNA framework.",
    "One of our goals is to enter Games Fleadh with the game.",
    "We have lead roles, but other group members can be of help as we all have experience with the area such as sounds design.",
    "We all plan on testing the product and giving feedback to the lead programmer."
)
$content = $newSlide.Shapes.Item(1)
$content.TextFrame.TextRange.Text = [string]::Join("`r", $bodyLines)

# ---------------------------------------------------------------
# 2. "Preferred Content" slide (slide 4) - add a trailing period to
#    the last three bullets, leaving the existing runs untouched.
# ---------------------------------------------------------------
$prefSlide = $p.Slides.Item(4)
$prefBody = $prefSlide.Shapes.Item(1).TextFrame.TextRange

$null = $prefBody.Characters(265, 84).InsertAfter(".")
$null = $prefBody.Characters(144, 120).InsertAfter(".")
$null = $prefBody.Characters(103, 40).InsertAfter(".")

# ---------------------------------------------------------------
# 3. Re-cache the "today" date placeholder text on the slide master
#    and every slide layout.
# ---------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $text) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $sh = $shapes.Item($k)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "06/02/2013"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes "06/02/2013"
}
